# "finish dev of vip effects"
# Remove the now-unused INT_recoveryCitizen column (column C) from the
# "dwelling" sheet; the remaining columns (poduction/power) shift left.
# Excel automatically re-indexes the shared string table (dropping the
# INT_recoveryCitizen entry) and updates every other sheet's references
# to the strings that shifted down one slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dwelling")

# Delete column C (INT_recoveryCitizen) entirely - shifts D/E left into C/D.
$ws.Columns.Item(3).Delete()

# Restore the active selection to what it becomes after the shift
# (old D4 selection -> now C4).
$ws.Range("C4").Select()
